$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LED PCB Assembly")

# Insert a new row at row 5 - shifts old rows 5..104 down to 6..105
$ws.Rows("5:5").Insert()

# ---- Row 3: Red SMT LEDs (surface mount) - Wurth Elektronik 150060RS75000 ----
$ws.Range("C3").Value = "W" + [char]0x00FC + "rth Elektronik"
$ws.Range("D3").Value = "150060RS75000"
$ws.Range("E3").HorizontalAlignment = -4108
$ws.Range("E3").Value = "DigiKey"
$ws.Range("F3").Value = "732-4978-1-ND"
$ws.Range("G3").HorizontalAlignment = -4108
$ws.Range("G3").Value = 56
$ws.Range("H3").NumberFormat = """$""#,##0.00"
$ws.Range("H3").HorizontalAlignment = -4108
$ws.Range("H3").Value = 0.142
$ws.Range("I3").Formula = "=G3*H3"
$ws.Range("J3").Value = "https://www.digikey.com/en/products/detail/w%C3%BCrth-elektronik/150060RS75000/4489901"

# ---- Row 4: Red SMT LEDs (surface mount, right angle) - Wurth Elektronik 155124RS73200 ----
$ws.Range("C4").Value = "W" + [char]0x00FC + "rth Elektronik"
$ws.Range("D4").Value = "155124RS73200"
$ws.Range("E4").HorizontalAlignment = -4108
$ws.Range("E4").Value = "DigiKey"
$ws.Range("F4").Value = "732-5025-1-ND"
$ws.Range("G4").HorizontalAlignment = -4108
$ws.Range("G4").Value = 2
$ws.Range("H4").NumberFormat = """$""#,##0.00"
$ws.Range("H4").HorizontalAlignment = -4108
$ws.Range("H4").Value = 0.2
$ws.Range("I4").Formula = "=G4*H4"
$ws.Range("J4").Value = "https://www.digikey.com/en/products/detail/w%C3%BCrth-elektronik/155124RS73200/4490041"

# ---- Row 5 (new row): Green SMT LEDs (surface mount) - Lite-On Inc. LTST-C191KGKT ----
$ws.Range("B5").Value = "Green SMT LEDs (surface mount) "
$ws.Range("C5").Value = "Lite-On Inc."
$ws.Range("D5").Value = "LTST-C191KGKT"
$ws.Range("E5").HorizontalAlignment = -4108
$ws.Range("E5").Value = "DigiKey"
$ws.Range("F5").Value = "160-1446-1-ND"
$ws.Range("I5").Formula = "=G5*H5"
$ws.Range("J5").Value = "https://www.digikey.com/en/products/detail/liteon/LTST-C191KGKT/386835"

# ---- Fix up the two real hyperlinks that shifted from J5/J7 to J6/J8 ----
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("J6"), "https://www.digikey.com/en/products/detail/analog-devices-inc-maxim-integrated/MAX6958AAEE-T/1521774")
$ws.Hyperlinks.Add($ws.Range("J8"), "https://www.digikey.com/en/products/detail/3m/D2510-6V0C-AR-WD/1886332")

# ---- Restore selection to match the author's final cursor position ----
$ws.Range("D17:D18").Select()
